$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quote row for 2025-11-12 (Excel serial date 45973), appended right
# after the last existing data row (row 68 -> A1:E68 becomes A1:E69).
$newRow = 69

$ws.Range("A$newRow").Value = 45973
# Match the date/number formatting already used by the column (copy the
# style from the row above so it keeps the same numFmt/date display).
$ws.Range("A$newRow").NumberFormat = $ws.Range("A68").NumberFormat

$ws.Range("B$newRow").Value = "22,0432"
$ws.Range("C$newRow").Value = "16,0723"
$ws.Range("D$newRow").Value = "15,5551"
$ws.Range("E$newRow").Value = "15,5551"
